$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$oldTimestamp = "2022-09-11 07:03:37"
$newTimestamp = "2022-09-11 20:59:53"

for ($row = 2; $row -le 64; $row++) {
    $cell = $ws.Cells.Item($row, 15)  # Column O is the 15th column
    if ($cell.Value2 -eq $oldTimestamp) {
        $cell.Value2 = $newTimestamp
    }
}
